# Update the "Metadata" sheet:
#  - translate the Publisher / Contact values from German to English
#  - fill in the previously-empty Description value

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"
$ws.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"
$ws.Range("B12").Value = "consent states - minimal subset WITHDRAWAL documents"
